# edit.ps1 — apply the "30. 11. 2021" wave update to ZBP_12_obavy_ztrata_prace.xlsx
#
# Sheet "data"   (index 1): new column AK — % values for the 30. 11. 2021 wave.
# Sheet "pocetR" (index 2): new column AJ — respondent counts for the same wave.
# Both sheets: extend the header date row, append the per-row figures, refresh
# the "aktualizace" (last-updated) date baked into the final footer row, and
# let Excel's own dimension tracking pick up the new A1:AK62 / A1:AJ25 extents.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "data"
$ws2 = $wb.Worksheets.Item(2)   # "pocetR"

# ---------------------------------------------------------------------------
# Sheet "data" — header cell AK1 (copy the look of AJ1, then set the text)
# ---------------------------------------------------------------------------
$ws1.Range("AJ1").Copy()
$ws1.Range("AK1").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("AK1").Value = "30. 11. 2021"

# ---------------------------------------------------------------------------
# Sheet "data" — AK2:AK61 data values (one per row, % as a fraction)
# ---------------------------------------------------------------------------
$ws1.Range("AK2").Value = 0.51
$ws1.Range("AK3").Value = 0.31
$ws1.Range("AK4").Value = 0.18
$ws1.Range("AK5").Value = 0.26
$ws1.Range("AK6").Value = 0.28
$ws1.Range("AK7").Value = 0.46
$ws1.Range("AK8").Value = 0.54
$ws1.Range("AK9").Value = 0.32
$ws1.Range("AK10").Value = 0.14
$ws1.Range("AK11").Value = 0.53
$ws1.Range("AK12").Value = 0.33
$ws1.Range("AK13").Value = 0.14
$ws1.Range("AK14").Value = 0.46
$ws1.Range("AK15").Value = 0.21
$ws1.Range("AK16").Value = 0.33
$ws1.Range("AK17").Value = 0.54
$ws1.Range("AK18").Value = 0.32
$ws1.Range("AK19").Value = 0.14
$ws1.Range("AK20").Value = 0.47
$ws1.Range("AK21").Value = 0.2
$ws1.Range("AK22").Value = 0.33
$ws1.Range("AK23").Value = 0.38
$ws1.Range("AK24").Value = 0.42
$ws1.Range("AK25").Value = 0.2
$ws1.Range("AK26").Value = 0.42
$ws1.Range("AK27").Value = 0.36
$ws1.Range("AK28").Value = 0.22
$ws1.Range("AK29").Value = 0.54
$ws1.Range("AK30").Value = 0.3
$ws1.Range("AK31").Value = 0.16
$ws1.Range("AK32").Value = 0.62
$ws1.Range("AK33").Value = 0.26
$ws1.Range("AK34").Value = 0.12
$ws1.Range("AK35").Value = 0.39
$ws1.Range("AK36").Value = 0.38
$ws1.Range("AK37").Value = 0.23
$ws1.Range("AK38").Value = 0.4
$ws1.Range("AK39").Value = 0.42
$ws1.Range("AK40").Value = 0.18
$ws1.Range("AK41").Value = 0.6
$ws1.Range("AK42").Value = 0.26
$ws1.Range("AK43").Value = 0.14
$ws1.Range("AK44").Value = 0.62
$ws1.Range("AK45").Value = 0.19
$ws1.Range("AK46").Value = 0.19
$ws1.Range("AK47").Value = 0.5
$ws1.Range("AK48").Value = 0.31
$ws1.Range("AK49").Value = 0.19
$ws1.Range("AK50").Value = 0.7
$ws1.Range("AK51").Value = 0.22
$ws1.Range("AK52").Value = 0.08
$ws1.Range("AK53").Value = 0.44
$ws1.Range("AK54").Value = 0.35
$ws1.Range("AK55").Value = 0.21
$ws1.Range("AK56").Value = 0.5600000000000001
$ws1.Range("AK57").Value = 0.32
$ws1.Range("AK58").Value = 0.12
$ws1.Range("AK59").Value = 0.63
$ws1.Range("AK60").Value = 0.17
$ws1.Range("AK61").Value = 0.2

# Footer label row 62 only carries a single cell (A62); bump the date inside
# the text, no AK62 cell exists in the source.
$ws1.Range("A62").Value = "Život během pandemie, Obavy ze ztráty práce, % respondentů celkově a ve skupinách, aktualizace 8. 12. 2021"

# ---------------------------------------------------------------------------
# Sheet "pocetR" — header cell AJ1 (copy the look of AI1, then set the text)
# ---------------------------------------------------------------------------
$ws2.Range("AI1").Copy()
$ws2.Range("AJ1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("AJ1").Value = "30. 11. 2021"

# ---------------------------------------------------------------------------
# Sheet "pocetR" — AJ2:AJ24 data values (respondent counts)
# ---------------------------------------------------------------------------
$ws2.Range("AJ2").Value = 1018
$ws2.Range("AJ3").Value = 93
$ws2.Range("AJ4").Value = 925
$ws2.Range("AJ5").Value = 784
$ws2.Range("AJ6").Value = 150
$ws2.Range("AJ7").Value = 10
$ws2.Range("AJ8").Value = 73
$ws2.Range("AJ9").Value = 751
$ws2.Range("AJ10").Value = 135
$ws2.Range("AJ11").Value = 70
$ws2.Range("AJ12").Value = 62
$ws2.Range("AJ13").Value = 374
$ws2.Range("AJ14").Value = 400
$ws2.Range("AJ15").Value = 244
$ws2.Range("AJ16").Value = 114
$ws2.Range("AJ17").Value = 304
$ws2.Range("AJ18").Value = 325
$ws2.Range("AJ19").Value = 157
$ws2.Range("AJ20").Value = 286
$ws2.Range("AJ21").Value = 102
$ws2.Range("AJ22").Value = 250
$ws2.Range("AJ23").Value = 143
$ws2.Range("AJ24").Value = 91

# Footer label row 25 carries the label in A25 plus blank placeholder cells
# across every other used column (B25:AI25); add the matching blank AJ25 so
# the row keeps a value in every column up to the new dimension, and bump the
# date inside the label text.
$ws2.Range("A25").Value = "Život během pandemie, Obavy ze ztráty práce, velikost dotázaného souboru celkově a ve skupinách, aktualizace 8. 12. 2021"
$ws2.Range("AJ25").Formula = "="""""
